# Update cryptos price/volume snapshot (GitHub Actions data refresh)
# D column holds price strings (possibly with a leading apostrophe to force
# text, since several values such as "1.00"/"7.60" would otherwise be
# auto-converted to numbers by Excel and lose their trailing zeros / the
# original text formatting used throughout the sheet).
# E column holds the 1h volume-change percentage strings, always kept as
# text because of the surrounding literal spaces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.850.52"
$ws.Range("E2").Value = "  +2.24%  "

$ws.Range("D3").Value = "3.471.32"
$ws.Range("E3").Value = "  +2.56%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'580.73"
$ws.Range("E5").Value = "  +1.01%  "

$ws.Range("D6").Value = "'147.37"
$ws.Range("E6").Value = "  +4.55%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.482"
$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("D9").Value = "'7.61"
$ws.Range("E9").Value = "  -0.80%  "

$ws.Range("E10").Value = "  +2.06%  "

$ws.Range("E11").Value = "  +4.07%  "

$ws.Range("D12").Value = "4.066.60"
$ws.Range("E12").Value = "  +2.63%  "

$ws.Range("D13").Value = "'29.68"
$ws.Range("E13").Value = "  +5.34%  "

$ws.Range("E14").Value = "  +2.64%  "

$ws.Range("D15").Value = "3.461.12"
$ws.Range("E15").Value = "  +2.81%  "

$ws.Range("E16").Value = "  +1.20%  "

$ws.Range("D17").Value = "62.877.67"
$ws.Range("E17").Value = "  +2.19%  "

$ws.Range("D18").Value = "'6.36"
$ws.Range("E18").Value = "  +4.21%  "

$ws.Range("D19").Value = "'14.35"
$ws.Range("E19").Value = "  +5.65%  "

$ws.Range("E20").Value = "  +2.80%  "

$ws.Range("D21").Value = "'388.40"
$ws.Range("E21").Value = "  -0.47%  "

$ws.Range("E22").Value = "  +2.69%  "

$ws.Range("D23").Value = "'74.72"
$ws.Range("E23").Value = "  -0.25%  "

$ws.Range("E24").Value = "  -0.19%  "

$ws.Range("D25").Value = "3.608.43"
$ws.Range("E25").Value = "  +2.48%  "

$ws.Range("E26").Value = "  +2.91%  "

$ws.Range("E27").Value = "  -9.00%  "

$ws.Range("D28").Value = "'7.60"
$ws.Range("E28").Value = "  +4.25%  "

$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  -0.27%  "

$ws.Range("D30").Value = "'8.17"
$ws.Range("E30").Value = "  +2.51%  "

$ws.Range("D31").Value = "'2.15"
$ws.Range("E31").Value = "  +0.46%  "

$ws.Range("E32").Value = "  +0.02%  "

$ws.Range("E33").Value = "  -0.46%  "

$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("D35").Value = "'5.30"
$ws.Range("E35").Value = "  +6.12%  "

$ws.Range("D36").Value = "'7.09"
$ws.Range("E36").Value = "  +3.16%  "

$ws.Range("D37").Value = "'1.60"
$ws.Range("E37").Value = "  +9.23%  "

$ws.Range("D38").Value = "'31.43"
$ws.Range("E38").Value = "  +20.66%  "

$ws.Range("D39").Value = "'170.36"
$ws.Range("E39").Value = "  +1.05%  "

$ws.Range("D40").Value = "3.510.77"
$ws.Range("E40").Value = "  +2.75%  "

$ws.Range("D41").Value = "'0.0768"
$ws.Range("E41").Value = "  +0.82%  "

$ws.Range("D42").Value = "'0.800"
$ws.Range("E42").Value = "  +2.91%  "

$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").Value = "'42.24"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("E45").Value = "  +4.78%  "

$ws.Range("E46").Value = "  +4.50%  "

$ws.Range("D47").Value = "2.612.05"
$ws.Range("E47").Value = "  +6.03%  "

$ws.Range("D48").Value = "'23.35"
$ws.Range("E48").Value = "  +3.41%  "

$ws.Range("D49").Value = "'6.75"
$ws.Range("E49").Value = "  +1.50%  "

$ws.Range("D50").Value = "'2.23"
$ws.Range("E50").Value = "  +10.93%  "

$ws.Range("E51").Value = "  +0.02%  "
